# Updated FIN_grids model - 2025-08-26 14:54
$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": reorder the comma-separated timeslice lists ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "S3aH2,S2aH2,S1aH2"
$wsEv.Range("C14").Value = "S2aH3,S3aH3,S3aH1,S2aH1,S1aH1,S1aH3"

# --- Sheet "re_profiles": re-shuffle the M/N rows 4-6 ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "S2"
$wsRe.Range("N4").Value = 0.59772867329870583
$wsRe.Range("M5").Value = "S3"
$wsRe.Range("N5").Value = 0.18498107227748917
$wsRe.Range("M6").Value = "S1"
$wsRe.Range("N6").Value = 0.41729025442380491

$wb.Save()
